$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 57: C57 formula changed (5876/36 -> 7618/36); D57 recomputes automatically ---
$ws.Range("C57").Formula = "=7618/36"
$ws.Range("D57").Formula = "=C57/(24*60)"

# --- Row 58: new data row for 四方坪站 (shared string index 2) on 2025-09-29 (serial 45929) ---
$ws.Range("A58").Value = 45929
$ws.Range("B58").Value = "四方坪站"
$ws.Range("C58").Formula = "=15852/127"
$ws.Range("D58").NumberFormat = "0.00%"
$ws.Range("D58").Formula = "=C58/(24*60)"
$ws.Range("E58").Formula = "=9170.63/127"
$ws.Range("F58").Formula = "=3207.95/127"
$ws.Range("G58").Formula = "=9170.63/(15852/60)"
$ws.Range("H58").Formula = "=373/127"

# --- Row 59: new data row for 高岭站 (shared string index 3) on 2025-09-29 (serial 45929) ---
$ws.Range("A59").Value = 45929
$ws.Range("B59").Value = "高岭站"
$ws.Range("C59").Formula = "=7057/36"
$ws.Range("D59").NumberFormat = "0.00%"
$ws.Range("D59").Formula = "=C59/(24*60)"
$ws.Range("E59").Formula = "=4957.18/36"
$ws.Range("F59").Formula = "=1241.19/36"
$ws.Range("G59").Formula = "=4957.18/(7057/60)"
$ws.Range("H59").Formula = "=189/36"

# --- Rows 60-63: the blank placeholder cell in column D is removed entirely ---
$ws.Range("D60").Clear()
$ws.Range("D61").Clear()
$ws.Range("D62").Clear()
$ws.Range("D63").Clear()

# --- Sheet view: scroll position and selected cell change ---
$win = $excel.ActiveWindow
$win.ScrollRow = 46
$win.ScrollColumn = 1
$ws.Range("I59").Select()
